$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.175.90'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '2.343.09'
$ws.Range("E3").Value = '  +1.68%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0786'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.72'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.66%  '
$ws.Range("E13").Value = '  +2.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("D15").Value = '2.705.53'
$ws.Range("E15").Value = '  +1.44%  '
$ws.Range("D16").Value = '2.333.40'
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.798'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.12%  '
$ws.Range("D18").Value = '43.098.65'
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.47%  '
$ws.Range("E21").Value = '  +0.41%  '
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.76'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.41'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  +0.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.44'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("E32").Value = '  +1.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0734'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.31'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.82%  '
$ws.Range("E35").Value = '  +4.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.37'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.07%  '
$ws.Range("E37").Value = '  -0.80%  '
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.76'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.22'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +19.30%  '
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '110.94'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -32.36%  '
$ws.Range("D43").Value = '1.938.16'
$ws.Range("E43").Value = '  -1.73%  '
$ws.Range("E44").Value = '  +1.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.74%  '
$ws.Range("E46").Value = '  +1.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.82%  '
$ws.Range("D48").Value = '2.570.30'
$ws.Range("E48").Value = '  +1.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("E50").Value = '  -3.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.21%  '
